# commit 28 - penambahan test suite collection dan scenario negative test case
#
# Adds 4 new worksheets (test-data sheets for new negative/edge scenarios),
# removes the now-unused "KOMPUTER" row from "data laptop", repositions the
# active-tab / first-visible-tab, and updates a handful of cell selections
# left over from interactive navigation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "data laptop" - drop the "KOMPUTER" row (row 4) that is no longer used
# ---------------------------------------------------------------------
$wsLaptop = $wb.Worksheets.Item("data laptop")
$wsLaptop.Range("A4").EntireRow.Delete()
$wsLaptop.Range("A7").Select()

# ---------------------------------------------------------------------
# 2) Cosmetic selection-only changes on existing sheets (left-over cursor
#    position from interactive editing)
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Activate()
$wsLogin.Range("F12").Select()

$wsRegister = $wb.Worksheets.Item("Register")
$wsRegister.Activate()
$wsRegister.Range("H10").Select()

$wsContactUs = $wb.Worksheets.Item("data contact us")
$wsContactUs.Activate()
$wsContactUs.Range("A1:B3").Select()

$wsSafepay = $wb.Worksheets.Item("checkout safepay")
$wsSafepay.Activate()
$wsSafepay.Range("A10").Select()

$wsNoLogin = $wb.Worksheets.Item("checkout wihout login")
$wsNoLogin.Activate()
$wsNoLogin.Range("E17").Select()

# ---------------------------------------------------------------------
# 3) New sheet: "login negative test case"
# ---------------------------------------------------------------------
$wsNegLogin = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNegLogin.Name = "login negative test case"
$wsNegLogin.Columns.Item(1).ColumnWidth = 25.5
$wsNegLogin.Columns.Item(2).ColumnWidth = 22

$wsNegLogin.Range("A1").Value = "username"
$wsNegLogin.Range("B1").Value = "password"
$wsNegLogin.Range("A2").Value = "yudhatesting"
$wsNegLogin.Range("B2").Value = 12345678
$wsNegLogin.Range("A3").Value = "yudhatesting2"
$wsNegLogin.Range("B3").Value = 12345678
$wsNegLogin.Range("A1:B3").Select()

# ---------------------------------------------------------------------
# 4) New sheet: "data gagal checkout"
# ---------------------------------------------------------------------
$wsGagalCheckout = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsGagalCheckout.Name = "data gagal checkout"
$wsGagalCheckout.Columns.Item(1).ColumnWidth = 24
$wsGagalCheckout.Columns.Item(2).ColumnWidth = 16.833333333333336

$wsGagalCheckout.Range("A1").Value = "username"
$wsGagalCheckout.Range("B1").Value = "password"
$wsGagalCheckout.Range("A2").Value = "yudhatesting"
$wsGagalCheckout.Range("B2").Value = 12345678
$wsGagalCheckout.Range("A3").Value = "yudhatesting2"
$wsGagalCheckout.Range("B3").Value = 12345678
$wsGagalCheckout.Range("G15").Select()

# ---------------------------------------------------------------------
# 5) New sheet: "data contact us gagal"
# ---------------------------------------------------------------------
$wsContactUsGagal = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsContactUsGagal.Name = "data contact us gagal"
$wsContactUsGagal.Columns.Item(1).ColumnWidth = 32.33333333333333
$wsContactUsGagal.Columns.Item(2).ColumnWidth = 31.666666666666668

$wsContactUsGagal.Range("A1").Value = "email"
$wsContactUsGagal.Range("B1").Value = "subject"

# The hyperlink's displayed/tooltip text intentionally differs from the cell's
# own text in the source workbook (sheet was cloned from "data contact us"
# and the cell text edited afterwards without touching the hyperlink) - so
# the TextToDisplay is set first and the cell text is overwritten afterwards.
$wsContactUsGagal.Range("A2").Value = "yudhatesting"
$wsContactUsGagal.Hyperlinks.Add($wsContactUsGagal.Range("A2"), "mailto:yudhatesting@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "yudhatesting@gmail.com")
$wsContactUsGagal.Range("A2").Style = "Hyperlink"
$wsContactUsGagal.Range("A2").Value = "yudhatesting"
$wsContactUsGagal.Range("B2").Value = "this product so cool!"

$wsContactUsGagal.Range("A3").Value = "yudhatesting2"
$wsContactUsGagal.Hyperlinks.Add($wsContactUsGagal.Range("A3"), "mailto:yudhatesting2@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "yudhatesting@")
$wsContactUsGagal.Range("A3").Style = "Hyperlink"
$wsContactUsGagal.Range("A3").Value = "yudhatesting2"
$wsContactUsGagal.Range("B3").Value = "I want to buy this product"

$wsContactUsGagal.Range("F8").Select()

# ---------------------------------------------------------------------
# 6) New sheet: "data keyword"
# ---------------------------------------------------------------------
$wsKeyword = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsKeyword.Name = "data keyword"

$wsKeyword.Range("A1").Value = "keyword"
$wsKeyword.Range("A2").Value = "naruto"
$wsKeyword.Range("A3").Value = "one piece"
$wsKeyword.Range("A4").Value = "bleach"
$wsKeyword.Range("L24").Select()

# ---------------------------------------------------------------------
# 7) Window view: active tab = "data contact us gagal" (10th sheet / index 9)
# ---------------------------------------------------------------------
$wsContactUsGagal.Activate()
$wsContactUsGagal.Range("F8").Select()
